$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: Terry Rozier -> Jrue Holiday (moved from row 19)
$ws.Range("A5").Value = "Jrue Holiday"
$ws.Range("B5").Value = "PG,SG"
$ws.Range("C5").Value = "Boston Celtics"

# Row 6: Paul George -> Caris LeVert (moved from row 14)
$ws.Range("A6").Value = "Caris LeVert"
$ws.Range("B6").Value = "SG,SF"
$ws.Range("C6").Value = "Cleveland Cavaliers"

# Row 7: Draymond Green -> Jayson Tatum (moved from row 16)
$ws.Range("A7").Value = "Jayson Tatum"
$ws.Range("B7").Value = "SF,PF"
$ws.Range("C7").Value = "Boston Celtics"

# Row 14: Caris LeVert -> Paul George (moved from row 6)
$ws.Range("A14").Value = "Paul George"
$ws.Range("B14").Value = "SG,SF,PF"
$ws.Range("C14").Value = "Philadelphia 76ers"

# Row 15: Anthony Edwards -> Draymond Green (moved from row 7)
$ws.Range("A15").Value = "Draymond Green"
$ws.Range("B15").Value = "PF,C"
$ws.Range("C15").Value = "Golden State Warriors"

# Row 16: Jayson Tatum -> Anthony Edwards (moved from row 15)
$ws.Range("A16").Value = "Anthony Edwards"
$ws.Range("B16").Value = "SG,SF"
$ws.Range("C16").Value = "Minnesota Timberwolves"

# Row 19: Jrue Holiday -> Anfernee Simons (new player, position unchanged)
$ws.Range("A19").Value = "Anfernee Simons"
$ws.Range("C19").Value = "Portland Trail Blazers"
